$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: DiemTB (average score) "8" -> "8,8" (keep stored as text, matching
# the original quote-prefixed text style used for this column)
$ws.Range("N3").Value = "'8,8"

# Row 4: DiemTB (average score) "9" -> "8,9"
$ws.Range("N4").Value = "'8,9"

# Row 3: DiemMonTO (Toán score) 5 -> "6,5" (now a text grade instead of a
# whole number)
$ws.Range("Z3").Value = "6,5"

# Row 4: DiemMonNV (Ngữ văn score) 5 -> "5,5"
$ws.Range("Y4").Value = "5,5"

# Update the active selection to match the author's final cursor position
$ws.Range("Y4").Select()
